$wb = $excel.ActiveWorkbook

# "Skill level" sheet: update the two form-control-linked cells to their new
# selections (dropdown "Choose Skill" -> 1, dropdown "Choose Personality 1" -> 3).
# Changing these cells drives the recalculation of every dependent formula
# (C3:F3, C5:F5, C10:F10, C12:F12) to match the committed values.
$wsSkillLevel = $wb.Worksheets.Item("Skill level")
$wsSkillLevel.Range("A3").Value = 1
$wsSkillLevel.Range("A5").Value = 3

# Keep the embedded combo-box form controls (the dropdowns driving A3/A5) in
# sync with the cells above.
$dropSkill = $wsSkillLevel.Shapes.Item("Vervolgkeuzelijst 4")
$dropSkill.ControlFormat.ListIndex = 1

$dropPersonality = $wsSkillLevel.Shapes.Item("Vervolgkeuzelijst 5")
$dropPersonality.ControlFormat.ListIndex = 3

# The active sheet moves from "Personality" back to "Skill level".
$wsSkillLevel.Activate()
